# Weekly update: insert a new Macroferia Regional de Talca - Alcachofa
# (Hortaliza) price record for a later week, pushing the old rows 48-50
# down to 49-51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 48 (old rows 48,49,50
# become 49,50,51).
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new record.
$ws.Cells.Item(48, 1).Value  = 5
$ws.Cells.Item(48, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(48, 3).Value  = "Maule"
$ws.Cells.Item(48, 4).Value  = 44461
$ws.Cells.Item(48, 5).Value  = 7
$ws.Cells.Item(48, 6).Value  = 100112013
$ws.Cells.Item(48, 7).Value  = "Alcachofa"
$ws.Cells.Item(48, 8).Value  = "Madrigal"
$ws.Cells.Item(48, 9).Value  = "Primera"
$ws.Cells.Item(48, 10).Value = 300
$ws.Cells.Item(48, 11).Value = 12000
$ws.Cells.Item(48, 12).Value = 12000
$ws.Cells.Item(48, 13).Value = 12000
$ws.Cells.Item(48, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(48, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(48, 16).Value = 300
$ws.Cells.Item(48, 17).Value = 40
$ws.Cells.Item(48, 18).Value = "Hortaliza"
